# This script applies the weekly order recalculation for Semana_14:
# several rows saw their "Unidades Pedido" (L) and "Diferencia Stock" (M)
# values increase, and the summary metrics "Total_Unidades" (C139) and
# "Total_Ajuste_Stock" (C150) are updated to reflect the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (Unidades Pedido, Diferencia Stock)
$updates = @{
    13  = @(4, 1)
    14  = @(5, 1)
    16  = @(4, 1)
    20  = @(6, 2)
    21  = @(8, 1)
    24  = @(11, 2)
    25  = @(7, 2)
    26  = @(11, 2)
    27  = @(6, 2)
    30  = @(16, 3)
    31  = @(8, 2)
    47  = @(2, 1)
    56  = @(9, 2)
    65  = @(4, 1)
    72  = @(2, 1)
    74  = @(2, 1)
    79  = @(5, 1)
    87  = @(7, 1)
    96  = @(9, 2)
    100 = @(2, 1)
    105 = @(2, 1)
    109 = @(5, 1)
    110 = @(9, 2)
    114 = @(2, 1)
    123 = @(2, 1)
    128 = @(3, 1)
    131 = @(2, 1)
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("L$row").Value = $values[0]
    $ws.Range("M$row").Value = $values[1]
}

# Update the summary metrics
$ws.Range("C139").Value = 372
$ws.Range("C150").Value = 38
